# Apply "add example values for arrayexpress-library" commit.

$wb = $excel.ActiveWorkbook

# 1. Bump version number 1.0.1 -> 1.0.2 on the isa_template sheet (cell B4).
$wsTemplate = $wb.Worksheets.Item("isa_template")
$wsTemplate.Range("B4").Value = "1.0.2"

# 2. Fill in example values in the "Library" sheet's data row (row 2).
$wsLibrary = $wb.Worksheets.Item("Library")

$wsLibrary.Range("F2").Value = "single-end"
$wsLibrary.Range("G2").Value = "DPBO"
$wsLibrary.Range("H2").Value = "http://purl.obolibrary.org/obo/DPBO_0000086"
$wsLibrary.Range("I2").Value = "Genome"
$wsLibrary.Range("J2").Value = "NCIT"
$wsLibrary.Range("K2").Value = "http://purl.obolibrary.org/obo/NCIT_C16629"
$wsLibrary.Range("L2").Value = "Whole Genome Sequencing"
$wsLibrary.Range("M2").Value = "NCIT"
$wsLibrary.Range("N2").Value = "http://purl.obolibrary.org/obo/NCIT_C101294"
$wsLibrary.Range("O2").Value = "Polymerase Chain Reaction"
$wsLibrary.Range("P2").Value = "NCIT"
$wsLibrary.Range("Q2").Value = "http://purl.obolibrary.org/obo/NCIT_C17003"
